$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Run 50" column (AZ) is removed entirely; this shifts the old "Mean"
# column (BA) left into AZ, matching Excel's native delete-column behaviour
# (dimension, spans and shared-string table all get recomputed on save).
$ws.Columns("AZ").Delete()

# Column header renamed from "Gen" to "MaxFES".
$ws.Range("A1").Value = "MaxFES"

# Column A values switch from generation counts to MaxFES fractions.
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# "Mean" column (now AZ after the delete) is recomputed over the
# remaining 50 runs (B:AY) instead of the original 51 (B:AZ).
$ws.Range("AZ2").Value = 466556548.3847353
$ws.Range("AZ3").Value = 243594544.4992435
$ws.Range("AZ4").Value = 30638668.12031402
$ws.Range("AZ5").Value = 1303078.54041253
$ws.Range("AZ6").Value = 409668.77896883
$ws.Range("AZ7").Value = 170651.35193839
$ws.Range("AZ8").Value = 101404.0333782
$ws.Range("AZ9").Value = 64353.00692167
$ws.Range("AZ10").Value = 47420.6148878
$ws.Range("AZ11").Value = 34257.2776684
$ws.Range("AZ12").Value = 28636.86392296
$ws.Range("AZ13").Value = 24378.62082797
$ws.Range("AZ14").Value = 22074.89990088
